$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E18").Value = 107

$ws.Range("E20").Value = 7

$ws.Range("F25").Value = 11
$ws.Range("H25").Value = 11

$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 5

$ws.Range("F28").Value = 12
$ws.Range("H28").Value = 12

$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 2

$ws.Range("E33").Value = 42

$ws.Range("E36").Value = 95

$ws.Range("F41").Value = 16
$ws.Range("H41").Value = 16

$ws.Range("E49").Value = 68

$ws.Range("E51").Value = 10

$ws.Range("E56").Value = 8

$ws.Range("E66").Value = 34

$ws.Range("E67").Value = 37
$ws.Range("F67").Value = 21
$ws.Range("H67").Value = 21

$ws.Range("F72").Value = 20
$ws.Range("H72").Value = 20

$ws.Range("E76").Value = 50

$ws.Range("E80").Value = 26

$ws.Range("E81").Value = 18
